$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 230, shifting existing rows 230:260 down to 231:261
$ws.Rows(230).Insert()

# Populate the newly inserted row 230 with the new weekly record
$ws.Range("A230").Value = 10
$ws.Range("B230").Value = "Vega Modelo de Temuco"
$ws.Range("C230").Value = "La Araucanía"
$ws.Range("D230").Value = 44776
$ws.Range("E230").Value = 9
$ws.Range("F230").Value = 100112043
$ws.Range("G230").Value = "Pepino dulce"
$ws.Range("H230").Value = "Cultivar IV Región"
$ws.Range("I230").Value = "Primera"
$ws.Range("J230").Value = 200
$ws.Range("K230").Value = 18000
$ws.Range("L230").Value = 19000
$ws.Range("M230").Value = 18500
$ws.Range("N230").Value = "$/bandeja 18 kilos"
$ws.Range("O230").Value = "Provincia de Limarí"
$ws.Range("P230").Value = 1028
$ws.Range("Q230").Value = 18
$ws.Range("R230").Value = "Hortaliza"
